# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1 (10:46 -> 11:16)
$ws.Range("A1").Value = "Datos actualizados a 23 de Marzo de 2020 a las 11:16"

# Update Madrid row (row 4): Casos activos, Recuperados, Muertes
$ws.Range("C4").Value = 2000
$ws.Range("D4").Value = 6681
$ws.Range("E4").Value = 1021

# Update Melilla row (row 53): Casos totales, Recuperados
$ws.Range("B53").Value = 28
$ws.Range("D53").Value = 28
